$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New meeting diary entry for 02/10/2023 (row 12), mirroring the formatting
# already used for the previous entry in row 10 (date / time-of-day / wrapped text).

# Date (column A) - reuse the existing date number format (m/d/yyyy)
$ws.Range("A12").Value = 45201
$ws.Range("A12").NumberFormat = "m/d/yy"

# Start / end time (columns B, C) - reuse the existing time number format (h:mm)
$ws.Range("B12").Value = 0.52083333333333337
$ws.Range("B12").NumberFormat = "h:mm"
$ws.Range("C12").Value = 0.5625
$ws.Range("C12").NumberFormat = "h:mm"

# Members present (column D)
$ws.Range("D12").Value = "All"

# Discussion notes (column E) - wrapped text like the other entries
$ws.Range("E12").Value = "Fixing expectations and deciding on how to proceed with the weather data"
$ws.Range("E12").WrapText = $true

# Match the row height used by the other meeting-entry rows
$ws.Rows.Item(12).RowHeight = 31.2

# Move the active selection, matching where the user left off after typing the entry
$ws.Range("A16").Select()
